$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 5935
$ws.Range("F5").Value = 69
$ws.Range("F6").Value = 2988
$ws.Range("F7").Value = 1277
$ws.Range("F11").Value = 725
$ws.Range("F12").Value = 275
$ws.Range("F13").Value = 4356
$ws.Range("F14").Value = 4356
$ws.Range("F17").Value = 121
$ws.Range("F20").Value = 69
$ws.Range("F21").Value = 6677
$ws.Range("F23").Value = 102
$ws.Range("F24").Value = 292
$ws.Range("F25").Value = 462
$ws.Range("G26").Value = 39.9
$ws.Range("F27").Value = 6255
$ws.Range("F30").Value = 1863
$ws.Range("F31").Value = 5988
$ws.Range("F34").Value = 98
$ws.Range("F36").Value = 420
$ws.Range("F37").Value = 4248
$ws.Range("F39").Value = 187
$ws.Range("F40").Value = 84
$ws.Range("F42").Value = 2
$ws.Range("F46").Value = 1005
$ws.Range("F48").Value = 332
$ws.Range("C49").Value = "【大会员提前抢】北京·第二届城市梦想动漫嘉年华（CDS）"
$ws.Range("F49").Value = 2064
$ws.Range("C50").Value = "【大会员提前抢】北京·万游引力国潮动漫嘉年华s7"
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F12").Value = 142
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1420
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1420
$ws.Range("F4").Value = 5935
$ws.Range("F5").Value = 2988
$ws.Range("F6").Value = 1277
$ws.Range("F12").Value = 275
$ws.Range("F13").Value = 4356
$ws.Range("F14").Value = 4356
$ws.Range("F17").Value = 121
$ws.Range("F20").Value = 69
$ws.Range("F21").Value = 6677
$ws.Range("F23").Value = 102
$ws.Range("F24").Value = 462
$ws.Range("G25").Value = 39.9
$ws.Range("F27").Value = 6255
$ws.Range("F29").Value = 1863
$ws.Range("F31").Value = 5988
$ws.Range("F35").Value = 98
$ws.Range("F37").Value = 420
$ws.Range("F38").Value = 4248
$ws.Range("F40").Value = 187
$ws.Range("F41").Value = 84
$ws.Range("F47").Value = 1005
$ws.Range("F49").Value = 332
$ws.Range("F50").Value = 142
